# Append a new sensor-reading row (row 56) to each of the four log
# worksheets, mirroring the structure of the existing rows (time, raw hex
# payload fields, and their decimal counterparts).

$wb = $excel.ActiveWorkbook

$rows = @{
    "ROW35-FE-LIFTER"  = @{
        A = "2025-03-06 15:42:06"
        B = "0x01,0x90 "
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x14,0x41,0x0c,"
        D = "0x01,0x90,"
        E = "0x d"
        F = 400
        G = "568631262647113770877196"
        H = 400
        I = 13
    }
    "ROW35-MID-LIFTER" = @{
        A = "2025-03-06 15:29:35"
        B = "0x01,0x90 "
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x15,0x41,0x0c,"
        D = "0x01,0x90,"
        E = "0x e"
        F = 400
        G = "568631262647113770942732"
        H = 400
        I = 14
    }
    "ROW02-FE-LIFTER"  = @{
        A = "2025-03-06 15:51:45"
        B = "0x01,0x90 "
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x06,0x41,0x0c,"
        D = "0x01,0x90,"
        E = "0xff"
        F = 400
        G = "568631262647113769959692"
        H = 400
        I = 255
    }
    "ROW02-MID-LIFTER" = @{
        A = "2025-03-06 15:41:15"
        B = "0x01,0x90 "
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,"
        D = "0x01,0x90,"
        E = "0x 3"
        F = 400
        G = "568631262647113769959692"
        H = 400
        I = 3
    }
}

foreach ($sheetName in $rows.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $data = $rows[$sheetName]
    $r = 56

    $ws.Cells.Item($r, 1).Value = $data.A
    $ws.Cells.Item($r, 2).Value = $data.B
    $ws.Cells.Item($r, 3).Value = $data.C
    $ws.Cells.Item($r, 4).Value = $data.D
    $ws.Cells.Item($r, 5).Value = $data.E
    $ws.Cells.Item($r, 6).Value = $data.F

    # Column G holds a long numeric-looking identifier that must stay text
    # (it exceeds double precision and would otherwise be coerced into a
    # rounded number). Prefixing with an apostrophe forces text entry, the
    # same way typing it directly into Excel would.
    $ws.Cells.Item($r, 7).Value = "'" + $data.G

    $ws.Cells.Item($r, 8).Value = $data.H
    $ws.Cells.Item($r, 9).Value = $data.I
}
